# Casos pendientes al 95%
# "Persona asignada" (columna C) se reasigna de "James Andres Urquiza"
# a "Luis Carlos Rincon Gordo" en las filas 2 y 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Luis Carlos Rincon Gordo"
$ws.Range("C3").Value = "Luis Carlos Rincon Gordo"

# Deja la selección donde quedó el usuario al guardar.
$ws.Range("M9").Select()
